$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old row 6 values (date, volume, unit label)
$ws.Range("D2").Value = 44208
$ws.Range("M2").Value = 210
$ws.Range("Q2").Value = "$/caja 14 kilos empedrada"

# Row 4 <- old row 2 values
$ws.Range("D4").Value = 44400
$ws.Range("M4").Value = 100
$ws.Range("Q4").Value = "$/caja 14 kilos"

# Row 5 <- old row 7 values
$ws.Range("D5").Value = 44491
$ws.Range("M5").Value = 180
$ws.Range("N5").Value = 9000
$ws.Range("O5").Value = 9000
$ws.Range("P5").Value = 9000
$ws.Range("S5").Value = 643

# Row 6 <- old row 4 values
$ws.Range("D6").Value = 44351
$ws.Range("M6").Value = 300

# Row 7 <- old row 5 values
$ws.Range("D7").Value = 44176
$ws.Range("M7").Value = 250
$ws.Range("N7").Value = 7000
$ws.Range("O7").Value = 7000
$ws.Range("P7").Value = 7000
$ws.Range("S7").Value = 500

# Row 8 <- old row 9 values
$ws.Range("D8").Value = 44162
$ws.Range("M8").Value = 120
$ws.Range("N8").Value = 7000
$ws.Range("O8").Value = 7000
$ws.Range("P8").Value = 7000
$ws.Range("Q8").Value = "$/caja 14 kilos empedrada"
$ws.Range("S8").Value = 500

# Row 9 <- old row 8 values
$ws.Range("D9").Value = 44397
$ws.Range("M9").Value = 60
$ws.Range("N9").Value = 11000
$ws.Range("O9").Value = 11000
$ws.Range("P9").Value = 11000
$ws.Range("Q9").Value = "$/caja 14 kilos"
$ws.Range("S9").Value = 786
